# Update "想去人数" (F) and "最低票价" (G) columns on the "展览" and
# "全部类型" sheets to match the latest scrape (commit 456a3b4).
#
# Sheet "展览" (rows keyed by row number on that sheet):
#   F2 642->641 ; G2 20 -> "不可售"
#   F4 53->56
#   F5 345->347
#   F6 423->427
#   F7 242->243
#   F8 13345->13388
#   F10 50->55
#   F11 5403->5424
#   F12 561->563
#   F13 29->30
#   F14 21->23
#   F15 45->46
#   F16 1211->1213
#   F18 149->151
#   F19 716->720
#   F20 2880->2882
#   F21 7390->7403
#   F23 3667->3670
#   F24 228->230
#   F25 56->57
#
# Sheet "全部类型" has the same events, offset by one extra row (row 6),
# so the same edits land on rows 2,4,5,7,8,9,11,12,13,14,15,16,17,19,20,21,23,25,26,27.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 641
$ws1.Range("G2").Value = "不可售"

$ws1.Range("F4").Value = 56
$ws1.Range("F5").Value = 347
$ws1.Range("F6").Value = 427
$ws1.Range("F7").Value = 243
$ws1.Range("F8").Value = 13388
$ws1.Range("F10").Value = 55
$ws1.Range("F11").Value = 5424
$ws1.Range("F12").Value = 563
$ws1.Range("F13").Value = 30
$ws1.Range("F14").Value = 23
$ws1.Range("F15").Value = 46
$ws1.Range("F16").Value = 1213
$ws1.Range("F18").Value = 151
$ws1.Range("F19").Value = 720
$ws1.Range("F20").Value = 2882
$ws1.Range("F21").Value = 7403
$ws1.Range("F23").Value = 3670
$ws1.Range("F24").Value = 230
$ws1.Range("F25").Value = 57

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 641
$ws4.Range("G2").Value = "不可售"

$ws4.Range("F4").Value = 56
$ws4.Range("F5").Value = 347
$ws4.Range("F7").Value = 427
$ws4.Range("F8").Value = 243
$ws4.Range("F9").Value = 13388
$ws4.Range("F11").Value = 55
$ws4.Range("F12").Value = 5424
$ws4.Range("F13").Value = 563
$ws4.Range("F14").Value = 30
$ws4.Range("F15").Value = 23
$ws4.Range("F16").Value = 46
$ws4.Range("F17").Value = 1213
$ws4.Range("F19").Value = 151
$ws4.Range("F20").Value = 720
$ws4.Range("F21").Value = 2882
$ws4.Range("F23").Value = 7403
$ws4.Range("F25").Value = 3670
$ws4.Range("F26").Value = 230
$ws4.Range("F27").Value = 57
